$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A394").Value = "2024-07-16 20:38:36"
$ws.Range("B394").Value = "Пользователь User 1 не удалось отправить сообщение"
$ws.Range("A395").Value = "2024-07-16 20:38:45"
$ws.Range("B395").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A396").Value = "2024-07-16 20:38:47"
$ws.Range("B396").Value = "Пользователь User 1 вышел из системы (logout)"
$ws.Range("A397").Value = "2024-07-16 20:38:55"
$ws.Range("B397").Value = "Пользователь User 2 вошёл в систему (log in)"
$ws.Range("A398").Value = "2024-07-16 20:38:56"
$ws.Range("B398").Value = "Пользователь User 2 вошёл в систему (cookie value)"
$ws.Range("A399").Value = "2024-07-16 20:39:03"
$ws.Range("B399").Value = "Пользователь User 2 написал сообщение 2: фкфукефкефке"
$ws.Range("A400").Value = "2024-07-17 22:59:12"
$ws.Range("B400").Value = "Пользователь User 1 вошёл в систему (log in)"
$ws.Range("A401").Value = "2024-07-17 22:59:14"
$ws.Range("B401").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A402").Value = "2024-07-17 22:59:20"
$ws.Range("B402").Value = "Пользователь User 1 написал сообщение 1: 12312313"
$ws.Range("A403").Value = "2024-07-17 22:59:44"
$ws.Range("B403").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A404").Value = "2024-07-17 23:02:22"
$ws.Range("B404").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A405").Value = "2024-07-17 23:33:52"
$ws.Range("B405").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A406").Value = "2024-07-17 23:34:03"
$ws.Range("B406").Value = "Пользователь User 1 не удалось отправить сообщение"
$ws.Range("A407").Value = "2024-07-17 23:34:31"
$ws.Range("B407").Value = "Пользователь User 1 не удалось отправить сообщение"
$ws.Range("A408").Value = "2024-07-17 23:36:14"
$ws.Range("B408").Value = "Пользователь User 1 не удалось отправить сообщение"
$ws.Range("A409").Value = "2024-07-17 23:36:29"
$ws.Range("B409").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A410").Value = "2024-07-17 23:55:16"
$ws.Range("B410").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A411").Value = "2024-07-17 23:55:24"
$ws.Range("B411").Value = "Пользователь User 1 написал сообщение 1:123123"
$ws.Range("A412").Value = "2024-07-17 23:56:54"
$ws.Range("B412").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A413").Value = "2024-07-17 23:57:08"
$ws.Range("B413").Value = "Пользователь User 1 написал сообщение 1: 22222"
$ws.Range("A414").Value = "2024-07-18 00:07:01"
$ws.Range("B414").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A415").Value = "2024-07-18 00:13:24"
$ws.Range("B415").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A416").Value = "2024-07-18 00:19:58"
$ws.Range("B416").Value = "Пользователь User 1 вошёл в систему (cookie value)"
$ws.Range("A417").Value = "2024-07-18 00:20:13"
$ws.Range("B417").Value = "Пользователь User 1 вышел из системы (logout)"
$ws.Range("A418").Value = "2024-07-18 00:20:20"
$ws.Range("B418").Value = "Пользователь User 2 вошёл в систему (log in)"
$ws.Range("A419").Value = "2024-07-18 00:20:22"
$ws.Range("B419").Value = "Пользователь User 2 вошёл в систему (cookie value)"
$ws.Range("A420").Value = "2024-07-18 00:20:33"
$ws.Range("B420").Value = "Пользователь User 2 написал сообщение 2: 2312323231"
